$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.031825702612842
$ws.Range("D2").Value = 1.041311234128588
$ws.Range("E2").Value = 1.031347056317505
$ws.Range("F2").Value = 1.0498088441016
$ws.Range("I2").Value = 1.037609870213082
$ws.Range("J2").Value = 1.03695929352739
$ws.Range("K2").Value = 1.044091108366749
$ws.Range("L2").Value = 1.034155437614137
$ws.Range("M2").Value = 1.05256487077546
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032826317500215
$ws.Range("D3").Value = 1.042099433132061
$ws.Range("E3").Value = 1.032198089628791
$ws.Range("F3").Value = 1.050730073968098
$ws.Range("I3").Value = 1.037832511416765
$ws.Range("J3").Value = 1.037601686063421
$ws.Range("K3").Value = 1.044689967995596
$ws.Range("L3").Value = 1.03481490173521
$ws.Range("M3").Value = 1.053298139988801
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033474039065305
$ws.Range("D4").Value = 1.042609314512821
$ws.Range("E4").Value = 1.032749369543917
$ws.Range("F4").Value = 1.051326352033497
$ws.Range("I4").Value = 1.037974874079927
$ws.Range("J4").Value = 1.03801702990991
$ws.Range("K4").Value = 1.045076679039764
$ws.Range("L4").Value = 1.035241582643736
$ws.Range("M4").Value = 1.053772165986827
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033746402167627
$ws.Range("D5").Value = 1.042823634939291
$ws.Range("E5").Value = 1.032981271293592
$ws.Range("F5").Value = 1.051577069487031
$ws.Range("I5").Value = 1.038034315530487
$ws.Range("J5").Value = 1.03819156142429
$ws.Range("K5").Value = 1.045239062078492
$ws.Range("L5").Value = 1.035420949772614
$ws.Range("M5").Value = 1.053971338238057
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033792136690165
$ws.Range("D6").Value = 1.042859618272524
$ws.Range("E6").Value = 1.033020217012735
$ws.Range("F6").Value = 1.051619168490188
$ws.Range("I6").Value = 1.038044272080601
$ws.Range("J6").Value = 1.038220861377255
$ws.Range("K6").Value = 1.045266315721369
$ws.Range("L6").Value = 1.035451065721164
$ws.Range("M6").Value = 1.054004773757535
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033477678155551
$ws.Range("D7").Value = 1.042612178406719
$ws.Range("E7").Value = 1.032752467663212
$ws.Range("F7").Value = 1.051329701968152
$ws.Range("I7").Value = 1.037975669942647
$ws.Range("J7").Value = 1.038019362320044
$ws.Range("K7").Value = 1.04507884955867
$ws.Range("L7").Value = 1.035243979394519
$ws.Range("M7").Value = 1.05377482776289
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.032163811604543
$ws.Range("D8").Value = 1.041577637171714
$ws.Range("E8").Value = 1.0316345412319
$ws.Range("F8").Value = 1.050120139955654
$ws.Range("I8").Value = 1.037685464728976
$ws.Range("J8").Value = 1.037176460375553
$ws.Range("K8").Value = 1.044293659077148
$ws.Range("L8").Value = 1.034378313677912
$ws.Range("M8").Value = 1.052812774760596
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029850600772033
$ws.Range("D9").Value = 1.03975365052736
$ws.Range("E9").Value = 1.029669286830514
$ws.Range("F9").Value = 1.047990170215742
$ws.Range("I9").Value = 1.037161082700744
$ws.Range("J9").Value = 1.035688686554783
$ws.Range("K9").Value = 1.042904031839743
$ws.Range("L9").Value = 1.032852655458056
$ws.Range("M9").Value = 1.051114124289293
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028309826647405
$ws.Range("D10").Value = 1.038537056918057
$ws.Range("E10").Value = 1.02836232008268
$ws.Range("F10").Value = 1.046571217865052
$ws.Range("I10").Value = 1.036802784554936
$ws.Range("J10").Value = 1.034695216552273
$ws.Range("K10").Value = 1.041973615309905
$ws.Range("L10").Value = 1.031835430989561
$ws.Range("M10").Value = 1.049979467088316
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027642983782997
$ws.Range("D11").Value = 1.03801012874233
$ws.Range("E11").Value = 1.027797160882846
$ws.Range("F11").Value = 1.045957052438766
$ws.Range("I11").Value = 1.036645578124839
$ws.Range("J11").Value = 1.034264656989911
$ws.Range("K11").Value = 1.041569796377041
$ws.Range("L11").Value = 1.031394942503397
$ws.Range("M11").Value = 1.04948763182629
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027395337433532
$ws.Range("D12").Value = 1.037814384860605
$ws.Range("E12").Value = 1.027587351611615
$ws.Range("F12").Value = 1.045728962632531
$ws.Range("I12").Value = 1.036586875492512
$ws.Range("J12").Value = 1.034104671444994
$ws.Range("K12").Value = 1.041419659201008
$ws.Range("L12").Value = 1.031231322611233
$ws.Range("M12").Value = 1.049304864843532
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027448456203374
$ws.Range("D13").Value = 1.037856373444814
$ws.Range("E13").Value = 1.027632351147447
$ws.Range("F13").Value = 1.045777886908879
$ws.Range("I13").Value = 1.036599481389203
$ws.Range("J13").Value = 1.034138991436784
$ws.Range("K13").Value = 1.041451870496484
$ws.Range("L13").Value = 1.031266419751742
$ws.Range("M13").Value = 1.049344072473903
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027622512261454
$ws.Range("D14").Value = 1.037993948879844
$ws.Range("E14").Value = 1.027779815610894
$ws.Range("F14").Value = 1.045938197671618
$ws.Range("I14").Value = 1.036640732051077
$ws.Range("J14").Value = 1.034251433684985
$ws.Range("K14").Value = 1.041557388860114
$ws.Range("L14").Value = 1.031381417686358
$ws.Range("M14").Value = 1.049472525825961
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027729760417333
$ws.Range("D15").Value = 1.038078711118992
$ws.Range("E15").Value = 1.027870688735297
$ws.Range("F15").Value = 1.046036975564579
$ws.Range("I15").Value = 1.036666106992093
$ws.Range("J15").Value = 1.034320705589655
$ws.Range("K15").Value = 1.041622383565506
$ws.Range("L15").Value = 1.031452271351395
$ws.Range("M15").Value = 1.049551659935468
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028354089634921
$ws.Range("D16").Value = 1.038572024644755
$ws.Range("E16").Value = 1.028399844056419
$ws.Range("F16").Value = 1.046611983358498
$ws.Range("I16").Value = 1.036813174446118
$ws.Range("J16").Value = 1.034723783408127
$ws.Range("K16").Value = 1.042000395650257
$ws.Range("L16").Value = 1.031864664321566
$ws.Range("M16").Value = 1.050012097665222
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028745801433266
$ws.Range("D17").Value = 1.038881431722435
$ws.Range("E17").Value = 1.028731974955534
$ws.Range("F17").Value = 1.046972738258356
$ws.Range("I17").Value = 1.036904874529285
$ws.Range("J17").Value = 1.034976522038083
$ws.Range("K17").Value = 1.042237260855183
$ws.Range("L17").Value = 1.032123341709372
$ws.Range("M17").Value = 1.050300779139543
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028974311392333
$ws.Range("D18").Value = 1.039061890532466
$ws.Range("E18").Value = 1.028925775137868
$ws.Range("F18").Value = 1.047183184573739
$ws.Range("I18").Value = 1.036958162751221
$ws.Range("J18").Value = 1.035123903451707
$ws.Range("K18").Value = 1.042375329340589
$ws.Range("L18").Value = 1.032274221532272
$ws.Range("M18").Value = 1.050469111814577
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.02905223264727
$ws.Range("D19").Value = 1.039123420070001
$ws.Range("E19").Value = 1.028991868489305
$ws.Range("F19").Value = 1.047254945384561
$ws.Range("I19").Value = 1.036976298909023
$ws.Range("J19").Value = 1.03517415045243
$ws.Range("K19").Value = 1.042422391671365
$ws.Range("L19").Value = 1.032325667238201
$ws.Range("M19").Value = 1.050526500341752
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028703771245717
$ws.Range("D20").Value = 1.038848236619252
$ws.Range("E20").Value = 1.028696332808068
$ws.Range("F20").Value = 1.046934030173341
$ws.Range("I20").Value = 1.036895056546368
$ws.Range("J20").Value = 1.034949409372058
$ws.Range("K20").Value = 1.042211856865523
$ws.Range("L20").Value = 1.032095588308808
$ws.Range("M20").Value = 1.050269811550682
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027571255741319
$ws.Range("D21").Value = 1.037953436867827
$ws.Range("E21").Value = 1.027736387791688
$ws.Range("F21").Value = 1.045890989074871
$ws.Range("I21").Value = 1.036628593289973
$ws.Range("J21").Value = 1.034218323801666
$ws.Range("K21").Value = 1.041526320214022
$ws.Range("L21").Value = 1.031347553734187
$ws.Range("M21").Value = 1.049434701647167
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026859480693635
$ws.Range("D22").Value = 1.037390729603626
$ws.Range("E22").Value = 1.027133503971444
$ws.Range("F22").Value = 1.045235411422279
$ws.Range("I22").Value = 1.036459268766154
$ws.Range("J22").Value = 1.033758334391626
$ws.Range("K22").Value = 1.041094481597219
$ws.Range("L22").Value = 1.030877218524364
$ws.Range("M22").Value = 1.048909187114683
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027236779235032
$ws.Range("D23").Value = 1.037689041572637
$ws.Range("E23").Value = 1.027453040110325
$ws.Range("F23").Value = 1.045582924013207
$ws.Range("I23").Value = 1.036549200271796
$ws.Range("J23").Value = 1.034002214263132
$ws.Range("K23").Value = 1.04132348437715
$ws.Range("L23").Value = 1.031126553423345
$ws.Range("M23").Value = 1.049187814478004
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028722762782149
$ws.Range("D24").Value = 1.03886323609805
$ws.Range("E24").Value = 1.028712437732181
$ws.Range("F24").Value = 1.046951520616401
$ws.Range("I24").Value = 1.036899493484767
$ws.Range("J24").Value = 1.034961660532346
$ws.Range("K24").Value = 1.042223336116081
$ws.Range("L24").Value = 1.032108128883081
$ws.Range("M24").Value = 1.050283804626682
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030448382467065
$ws.Range("D25").Value = 1.040225305730964
$ws.Range("E25").Value = 1.030176792040691
$ws.Range("F25").Value = 1.048540642014827
$ws.Range("I25").Value = 1.03729818561901
$ws.Range("J25").Value = 1.036073600478375
$ws.Range("K25").Value = 1.043263992174561
$ws.Range("L25").Value = 1.033247098970678
$ws.Range("M25").Value = 1.051553662137748
